$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:H1 matching the style of existing header (A1:E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header style from E1 (which already has the bold/bordered header style) to F1:H1
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Boolean values for KNN_Outliers_MAD (F), SVM_Outliers_MAD (G), RF_Outliers_MAD (H)
# All FALSE in F and G; H is TRUE only for rows 3 and 21
$trueRows = @(3, 21)

for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 6).Value = $false   # F column - KNN_Outliers_MAD
    $ws.Cells.Item($row, 7).Value = $false   # G column - SVM_Outliers_MAD
    if ($trueRows -contains $row) {
        $ws.Cells.Item($row, 8).Value = $true
    } else {
        $ws.Cells.Item($row, 8).Value = $false
    }
}
